$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2 (Execute signature:
# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
# MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)

# 1) Tag name: RSDDRId -> RSDRId
$d.Content.Find.Execute("RSDDRId", $false, $false, $false, $false, $false, `
    $true, 1, $false, "RSDRId", 2)

# 2) Corresponding field label: ID DDR partagé -> ID Demande partagé
$d.Content.Find.Execute("ID DDR partagé", $false, $false, $false, $false, $false, `
    $true, 1, $false, "ID Demande partagé", 2)

# 3) Description: drop the trailing clause and append a value-pattern line,
#    separated by a manual line break (^l -> <w:br/>) inside the same run.
$d.Content.Find.Execute( `
    "Identifiant unique partagé de la demande de ressource à laquelle l'expéditeur répond", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Identifiant unique partagé de la demande de ressource^l{orgID}.D.{ID unique de la demande dans le système émetteur}", `
    2)

# 4) Fill the previously-empty "Exemple" cell for the RSDRId row
#    (row 3, column 6 of the first table) with the sample value.
$table = $d.Tables.Item(1)
$cell = $table.Cell(3, 6)
$cell.Range.Text = "fr.health.samu770.D.1249875"

# 5) Field label: Réponse à la demande de concours -> ... de ressources
$d.Content.Find.Execute("Réponse à la demande de concours", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Réponse à la demande de ressources", 2)

# 6) ENUM values for the "answer" field
$d.Content.Find.Execute("(ENUM : OUI, NON, PARTIEL, DIFFERE)", $false, $false, $false, `
    $false, $false, $true, 1, $false, "(ENUM : ACCEPTEE, REFUSEE, PARTIELLE, DIFFEREE)", 2)
